$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date-serial value (45182 -> 45184) for every
# data row (rows 2 through 483). Update them all in one shot via a Range.
$ws.Range("C2:C483").Value = 45184
